# Update ticket-sales / view counts (column F) on both the "展览" sheet
# and its duplicate on the "全部类型" sheet, per the regenerated
# gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13550
$ws1.Range("F5").Value  = 1030
$ws1.Range("F6").Value  = 21
$ws1.Range("F7").Value  = 1742
$ws1.Range("F8").Value  = 143
$ws1.Range("F10").Value = 81
$ws1.Range("F13").Value = 13563
$ws1.Range("F14").Value = 340
$ws1.Range("F16").Value = 8966
$ws1.Range("F18").Value = 8059
$ws1.Range("F27").Value = 15
$ws1.Range("F31").Value = 186

# --- Sheet "全部类型" (all types, same rows as 展览 but shifted by +2) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 13550
$ws4.Range("F5").Value  = 1030
$ws4.Range("F6").Value  = 21
$ws4.Range("F7").Value  = 1742
$ws4.Range("F8").Value  = 143
$ws4.Range("F10").Value = 81
$ws4.Range("F13").Value = 13563
$ws4.Range("F14").Value = 340
$ws4.Range("F16").Value = 8966
$ws4.Range("F18").Value = 8059
$ws4.Range("F27").Value = 15
$ws4.Range("F33").Value = 186
